$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corporate/person name separator fixes (comma -> period)
$textFixes = @(
  @("E34", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("F34", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("E39", "RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH"),
  @("E40", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
  @("E43", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("F43", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("E63", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
  @("E99", "FERNANDEZ. MARIO HUGO"),
  @("E101", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("F101", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
  @("E102", "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"),
  @("E124", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
  @("E133", "TRABICHET MARIA. VERGARA ADEL Y OTRA"),
  @("F133", "TRABICHET MARIA. VERGARA ADEL Y OTRA"),
  @("E141", "RICCOTTI. MARIANA EDITH"),
  @("F147", "MERCANZINI. GASTON ARIEL"),
  @("E173", "FERNANDEZ. MARIO HUGO"),
  @("E196", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH")
)

foreach ($p in $textFixes) {
  $ws.Range($p[0]).Value = $p[1]
}

# Numeric-looking text values re-scraped without thousands separators,
# using comma-as-decimal -> dot-as-decimal. Values must remain TEXT
# (quote-prefixed) since the source column stores these as strings.
$numberFixes = @(
  @("H2", "'3774.00"),
  @("H3", "'43490.00"),
  @("H4", "'458000.00"),
  @("H5", "'42531.50"),
  @("H6", "'120.00"),
  @("H7", "'17.45"),
  @("H8", "'2036.42"),
  @("H9", "'121000.00"),
  @("H10", "'2999.99"),
  @("H11", "'31165.08"),
  @("H12", "'87334.15"),
  @("H13", "'29750.00"),
  @("H14", "'39478.94"),
  @("H15", "'1050.00"),
  @("H16", "'5955.70"),
  @("H17", "'9624.44"),
  @("H18", "'493.00"),
  @("H19", "'13630.75"),
  @("H20", "'9810.00"),
  @("H21", "'160.00"),
  @("H22", "'179.00"),
  @("H23", "'797.00"),
  @("H24", "'238.90"),
  @("H25", "'529.22"),
  @("H26", "'84.99"),
  @("H27", "'51.90"),
  @("H28", "'256.56"),
  @("H29", "'30.00"),
  @("H30", "'57290.00"),
  @("H31", "'42338.32"),
  @("H32", "'275.00"),
  @("H33", "'2384.58"),
  @("H34", "'28.80"),
  @("H35", "'25.04"),
  @("H36", "'2297.09"),
  @("H37", "'756.00"),
  @("H38", "'21024.50"),
  @("H39", "'39627.25"),
  @("H40", "'280.00"),
  @("H41", "'16177.15"),
  @("H42", "'801.00"),
  @("H43", "'295.22"),
  @("H44", "'8778.76"),
  @("H45", "'618.00"),
  @("H46", "'980.00"),
  @("H47", "'1394.49"),
  @("H48", "'1100.00"),
  @("H49", "'41.60"),
  @("H50", "'830.00"),
  @("H51", "'96.00"),
  @("H52", "'1397.20"),
  @("H53", "'287.84"),
  @("H54", "'23.63"),
  @("H55", "'45.21"),
  @("H56", "'2500.00"),
  @("H57", "'115.12"),
  @("H58", "'13.30"),
  @("H59", "'1936.10"),
  @("H60", "'800.00"),
  @("H61", "'208.35"),
  @("H62", "'65.00"),
  @("H63", "'160.00"),
  @("H64", "'899.64"),
  @("H65", "'40.48"),
  @("H66", "'380.00"),
  @("H67", "'3750.00"),
  @("H68", "'92.35"),
  @("H69", "'3035.00"),
  @("H70", "'1920.00"),
  @("H71", "'240.00"),
  @("H72", "'1134.00"),
  @("H73", "'6924.96"),
  @("H74", "'80.40"),
  @("H75", "'2881.40"),
  @("H76", "'1201.71"),
  @("H77", "'746.25"),
  @("H78", "'109.70"),
  @("H79", "'335.00"),
  @("H80", "'311.88"),
  @("H81", "'72.00"),
  @("H82", "'320.00"),
  @("H83", "'8288.00"),
  @("H84", "'7820.66"),
  @("H85", "'960.00"),
  @("H86", "'302.50"),
  @("H87", "'1355.00"),
  @("H88", "'3800.00"),
  @("H89", "'51312.00"),
  @("H90", "'320.00"),
  @("H91", "'8360.00"),
  @("H92", "'34437.00"),
  @("H93", "'11000.00"),
  @("H94", "'98090.00"),
  @("H95", "'4574.70"),
  @("H96", "'1485.00"),
  @("H97", "'44.00"),
  @("H98", "'899.00"),
  @("H99", "'367.50"),
  @("H100", "'1508.00"),
  @("H101", "'246.80"),
  @("H102", "'850.00"),
  @("H103", "'57.95"),
  @("H104", "'110.00"),
  @("H105", "'78.00"),
  @("H106", "'42.00"),
  @("H107", "'282.00"),
  @("H108", "'20.00"),
  @("H109", "'150.00"),
  @("H110", "'10000.00"),
  @("H111", "'5000.00"),
  @("H112", "'62000.00"),
  @("H113", "'222570.00"),
  @("H114", "'4650.60"),
  @("H115", "'3.79"),
  @("H116", "'19.20"),
  @("H117", "'49.68"),
  @("H118", "'15654.25"),
  @("H119", "'178.92"),
  @("H120", "'18.80"),
  @("H121", "'1522.64"),
  @("H122", "'21.45"),
  @("H123", "'2896.00"),
  @("H124", "'460.00"),
  @("H125", "'305.00"),
  @("H126", "'31.40"),
  @("H127", "'34.00"),
  @("H128", "'1415.70"),
  @("H129", "'102.00"),
  @("H130", "'216.00"),
  @("H131", "'1648.75"),
  @("H132", "'3237.00"),
  @("H133", "'20587.59"),
  @("H134", "'4571.15"),
  @("H135", "'112.00"),
  @("H136", "'472.75"),
  @("H137", "'133.45"),
  @("H138", "'422.30"),
  @("H139", "'328.67"),
  @("H140", "'37.03"),
  @("H141", "'1000.00"),
  @("H142", "'95397.74"),
  @("H143", "'14940.00"),
  @("H144", "'15628.00"),
  @("H145", "'1500.00"),
  @("H146", "'1060.00"),
  @("H147", "'6000.00"),
  @("H148", "'2000.00"),
  @("H149", "'4782.69"),
  @("H150", "'1029.60"),
  @("H151", "'795.15"),
  @("H152", "'406.40"),
  @("H153", "'1637.00"),
  @("H154", "'1727.00"),
  @("H155", "'131505.20"),
  @("H156", "'7320.00"),
  @("H157", "'1600.00"),
  @("H158", "'1000.00"),
  @("H159", "'4686.34"),
  @("H160", "'384.00"),
  @("H161", "'800.00"),
  @("H162", "'1000.00"),
  @("H163", "'12577.50"),
  @("H164", "'1500.00"),
  @("H165", "'950.00"),
  @("H166", "'750.00"),
  @("H167", "'2850.00"),
  @("H168", "'1500.00"),
  @("H169", "'200.00"),
  @("H170", "'290.00"),
  @("H171", "'16080.00"),
  @("H172", "'1300.00"),
  @("H173", "'120.00"),
  @("H174", "'484.00"),
  @("H175", "'135.00"),
  @("H176", "'42955.00"),
  @("H177", "'3680.00"),
  @("H178", "'12336.50"),
  @("H179", "'285.95"),
  @("H180", "'550.00"),
  @("H181", "'1197.30"),
  @("H182", "'5352.00"),
  @("H183", "'1125.00"),
  @("H184", "'2014.80"),
  @("H185", "'9700.00"),
  @("H186", "'4155.50"),
  @("H187", "'1285.42"),
  @("H188", "'150.00"),
  @("H189", "'637.66"),
  @("H190", "'35.92"),
  @("H191", "'2.60"),
  @("H192", "'3527.80"),
  @("H193", "'210.00"),
  @("H194", "'4755.32"),
  @("H195", "'140.00"),
  @("H196", "'1825.00"),
  @("H197", "'237.00"),
  @("H198", "'279.64"),
  @("H199", "'780.00"),
  @("H200", "'309.20"),
  @("H201", "'50.00"),
  @("H202", "'1558.50"),
  @("H203", "'120.40"),
  @("H204", "'16431.00"),
  @("H205", "'12905.50"),
  @("H206", "'3256.00"),
  @("H207", "'470.00"),
  @("H208", "'475.00"),
  @("H209", "'1726.80"),
  @("H210", "'4144.99"),
  @("H211", "'7813.49"),
  @("H212", "'17920.00"),
  @("H213", "'690471.47"),
  @("H214", "'890.00"),
  @("H215", "'216500.00"),
  @("H216", "'52000.00"),
  @("H217", "'25905.75"),
  @("H218", "'61000.00"),
  @("H219", "'20000.00"),
  @("H220", "'89500.00"),
  @("H221", "'145260.00"),
  @("H222", "'12000.00"),
  @("H223", "'223000.00"),
  @("H224", "'111500.00"),
  @("H225", "'223000.00"),
  @("H226", "'5950.00"),
  @("H227", "'617607.56"),
  @("H228", "'502317.08"),
  @("H229", "'11077.12"),
  @("H230", "'62400.00"),
  @("H231", "'1053000.00"),
  @("H232", "'6000.00"),
  @("H233", "'1855.00")
)

foreach ($p in $numberFixes) {
  $ws.Range($p[0]).Value = $p[1]
}

Write-Host "Applied $($textFixes.Count) text fixes and $($numberFixes.Count) number fixes"
